$wb = $excel.ActiveWorkbook

# zh-cn sheet: Priority -> "ht" and Latest Handoff Datetime refreshed for the
# "0de7f9ee..." handoff batch (rows 4-7)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4:E7").Value = "ht"
$wsZh.Range("H4:H7").Value = "2016-08-27 12:29:24"

# de-de sheet: same Priority refresh plus the handoff generation timestamp bump
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4:E7").Value = "ht"
$wsDe.Range("H4:H7").Value = "2016-08-27 12:29:29"

# Overview sheet: "Latest HO Xliff Generate Date" shares the same underlying
# value as de-de's refreshed handoff timestamp, so it picks up the same bump
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-27 12:29:29"
